$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.358.37"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "1.824.86"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'315.21"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D7").Value = "'0.4475"
$ws.Range("E7").Value = "  -2.26%  "
$ws.Range("D8").Value = "'0.3776"
$ws.Range("E8").Value = "  +0.95%  "
$ws.Range("D9").Value = "'0.07472"
$ws.Range("E9").Value = "  +1.94%  "
$ws.Range("D10").Value = "'0.8862"
$ws.Range("E10").Value = "  +2.90%  "
$ws.Range("D11").Value = "'21.02"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").Value = "1.825.05"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "'6.752"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").Value = "'5.451"
$ws.Range("E14").Value = "  +1.80%  "
$ws.Range("D15").Value = "'93.78"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "'0.07121"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "'0.000008798"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "'15.16"
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("D21").Value = "27.388.10"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("E22").Value = "  +4.08%  "
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("D24").Value = "'1.964"
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("D25").Value = "'151.54"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "'2.308"
$ws.Range("E26").Value = "  +4.06%  "
$ws.Range("D27").Value = "'18.72"
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("D28").Value = "'5.413"
$ws.Range("E28").Value = "  +2.72%  "
$ws.Range("D29").Value = "'117.81"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("D31").Value = "'0.7910"
$ws.Range("E31").Value = "  +2.24%  "
$ws.Range("D32").Value = "'1.206"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("D33").Value = "'4.602"
$ws.Range("E33").Value = "  +2.92%  "
$ws.Range("D34").Value = "'2.922"
$ws.Range("E34").Value = "  -1.70%  "
$ws.Range("D35").Value = "'1.0000"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'1.110"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").Value = "'0.01990"
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("D38").Value = "'0.05307"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").Value = "'7.319"
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("D41").Value = "'2.874"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").Value = "'0.1720"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "'2.313"
$ws.Range("E43").Value = "  +15.92%  "
$ws.Range("D44").Value = "'8.676"
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("D45").Value = "'0.5116"
$ws.Range("E45").Value = "  -2.93%  "
$ws.Range("D46").Value = "'10.68"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").Value = "'1.697"
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("D48").Value = "'105.44"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "'0.06412"
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("D51").Value = "'66.01"
$ws.Range("E51").Value = "  +3.86%  "
